$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N (14th column),
# shifting existing N:P data into O:Q
$ws.Columns("N").Insert()

# Give the newly inserted column N the same width as column M
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selection on this sheet to reflect where the user clicked after the edit
$ws.Range("R8").Select()

# Make this sheet the active one (moves tabSelected here from "Floating Interest Rates")
$ws.Activate()
